$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the dataset. It belongs right
# above the existing row 140 ("Inferno" / "Primera" / Región de Arica y
# Parinacota, $/caja 12 kilos), so insert a fresh row there and push all
# the following rows (old 140..181) down to (141..182).
$ws.Rows(140).Insert()

# Populate the newly inserted row 140 with the new record's data.
$ws.Range("A140").Value = 9
$ws.Range("B140").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C140").Value = "Metropolitana"
$ws.Range("D140").Value = 44508
$ws.Range("E140").Value = 13
$ws.Range("F140").Value = 100112021
$ws.Range("G140").Value = "Ají"
$ws.Range("H140").Value = "Inferno"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 52
$ws.Range("K140").Value = 19000
$ws.Range("L140").Value = 20000
$ws.Range("M140").Value = 19500
$ws.Range("N140").Value = "$/caja 12 kilos"
$ws.Range("O140").Value = "Región de Arica y Parinacota"
$ws.Range("P140").Value = 1625
$ws.Range("Q140").Value = 12
$ws.Range("R140").Value = "Hortaliza"
